# 8.3.1.2.xlsx — add the "2020" column (K) to the small/medium enterprise
# employment-share table and move the active selection, per the target diff.
#
# Column K picks up a new medium-weight top/bottom border, two new number
# formats/fonts, etc. Rather than re-deriving every style property from
# scratch (and risking the host mistakenly inventing brand-new font/border
# table entries for combinations that already exist), we seed each new cell
# by copying the *closest* already-formatted neighbour's format onto it and
# then nudging only the handful of properties that actually differ. That
# mirrors how this row/column was almost certainly built in real Excel
# (copy the adjacent year column, retype the header/values, restyle the
# numeric format/font/border for the new "2020" column) and keeps the
# workbook's shared style tables from growing more than necessary.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# K3: header year "2020" — same look as the I3 header cell (bold 10pt
# Times New Roman, medium top+bottom border), just not centered.
# ---------------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("K3").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("K3").Font.Size = 10
$ws.Range("K3").VerticalAlignment = -4107  # xlBottom (clear inherited vcenter)
$ws.Range("K3").Value = 2020

# ---------------------------------------------------------------------
# K4: "Small enterprises" 2020 value — same base look as J4 (9pt Times
# New Roman), but in the new Kyrghyz Times font, right aligned, with a
# thousands number format and a new medium top border marking the start
# of the 2020 column block.
# ---------------------------------------------------------------------
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("K4").Font.Name = "Kyrghyz Times"
$ws.Range("K4").NumberFormat = "#,##0.0"
$ws.Range("K4").HorizontalAlignment = -4152  # xlRight
$ws.Range("K4").VerticalAlignment = -4107    # xlBottom (clear inherited vcenter)
$ws.Range("K4").Borders.Item(8).Weight = -4138  # xlEdgeTop, xlMedium
$ws.Range("K4").Value = 2.8218550629805335

# ---------------------------------------------------------------------
# K5: "Medium-sized enterprises" 2020 value — same base look as J5
# (9pt Times New Roman, medium bottom border already present), switched
# to Kyrghyz Times, right aligned, thousands number format.
# ---------------------------------------------------------------------
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("K5").Font.Name = "Kyrghyz Times"
$ws.Range("K5").NumberFormat = "#,##0.0"
$ws.Range("K5").HorizontalAlignment = -4152  # xlRight
$ws.Range("K5").VerticalAlignment = -4107    # xlBottom (clear inherited vcenter)
$ws.Range("K5").Value = 1.3005071159823327

# ---------------------------------------------------------------------
# Move the active selection (matches the saved sheet view in the diff).
# ---------------------------------------------------------------------
$ws.Range("L8").Select()
